$wb = $excel.ActiveWorkbook

# sheet1 (Worksheets.Item(1))
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 5072
$ws.Range("F5").Value = 5072
$ws.Range("F7").Value = 150
$ws.Range("F10").Value = 212
$ws.Range("F12").Value = 8381
$ws.Range("F13").Value = 8381
$ws.Range("F14").Value = 273
$ws.Range("F15").Value = 136
$ws.Range("F17").Value = 609
$ws.Range("F18").Value = 2524
$ws.Range("F19").Value = 436
$ws.Range("F21").Value = 2300
$ws.Range("G21").Value = 63
$ws.Range("F23").Value = 27
$ws.Range("F24").Value = 2522
$ws.Range("F27").Value = 6400
$ws.Range("F28").Value = 183
$ws.Range("F29").Value = 68
$ws.Range("F32").Value = 458
$ws.Range("F33").Value = 6890
$ws.Range("F40").Value = 25
$ws.Range("F43").Value = 2528
$ws.Range("F48").Value = 516
$ws.Range("F49").Value = 2218
$ws.Range("F50").Value = 74

# sheet2 (Worksheets.Item(2))
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 166
$ws.Range("F5").Value = 55
$ws.Range("F6").Value = 14
$ws.Range("F11").Value = 39

# sheet4 (Worksheets.Item(4))
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 5072
$ws.Range("F4").Value = 5072
$ws.Range("F6").Value = 150
$ws.Range("F9").Value = 212
$ws.Range("F11").Value = 8381
$ws.Range("F12").Value = 8381
$ws.Range("F13").Value = 273
$ws.Range("F14").Value = 136
$ws.Range("F15").Value = 609
$ws.Range("F16").Value = 2524
$ws.Range("F17").Value = 166
$ws.Range("F19").Value = 2300
$ws.Range("G19").Value = 63
$ws.Range("F20").Value = 55
$ws.Range("F22").Value = 2522
$ws.Range("F24").Value = 14
$ws.Range("F27").Value = 6400
$ws.Range("F28").Value = 183
$ws.Range("F29").Value = 68
$ws.Range("F32").Value = 458
$ws.Range("F33").Value = 6890
$ws.Range("F39").Value = 25
$ws.Range("F40").Value = 39
$ws.Range("F42").Value = 2528
$ws.Range("F47").Value = 516
$ws.Range("F49").Value = 2218
$ws.Range("F50").Value = 74
